$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Helper: clone the "data row" cell formatting (style index used by the other
# populated rows, e.g. row 263) onto a target row, then fill in the values.
# Using Copy/PasteSpecial(xlPasteFormats) instead of touching Font/Borders by
# hand lets the engine de-duplicate back onto the existing style record
# instead of fabricating a brand-new (duplicate) style entry.
# ---------------------------------------------------------------------------

function Set-EventRow {
    param($row, $date, $eventText, $locationText, $cityText, $linkUrl)

    # Clone formatting from the template row (263) which already carries the
    # "populated data row" style for columns B:E.
    $ws.Range("B263:E263").Copy() | Out-Null
    $ws.Range("B" + $row + ":E" + $row).PasteSpecial(-4122) | Out-Null

    $ws.Range("A" + $row).Value = $date
    $ws.Range("B" + $row).Value = $eventText
    $ws.Range("C" + $row).Value = $locationText
    $ws.Range("D" + $row).Value = $cityText
    $ws.Range("E" + $row).Value = $linkUrl

    # Give the URL text the same "hyperlink look" (single underline, indexed
    # color 11) used by the other link cells in the sheet. Splitting the run
    # into two adjoining Characters() calls keeps the final shared-string as
    # ONE rich run spanning the whole text instead of collapsing into a
    # plain (non-rich) string tied to the cell-level style.
    $len = $linkUrl.Length
    $firstPart = $ws.Range("E" + $row).Characters(1, $len - 1)
    $firstPart.Font.Underline = 2
    $firstPart.Font.ColorIndex = 4
    $lastPart = $ws.Range("E" + $row).Characters($len, 1)
    $lastPart.Font.Underline = 2
    $lastPart.Font.ColorIndex = 4

    # Register the real clickable hyperlink (creates the relationship +
    # worksheet hyperlink entry).
    $ws.Hyperlinks.Add($ws.Range("E" + $row), $linkUrl, "", "", $linkUrl) | Out-Null

    # Hyperlinks.Add forces the built-in "Hyperlink" cell style onto E<row>;
    # reapply the template formatting so the cell keeps using the same style
    # as the rest of the row (matches the original workbook's convention of
    # storing link styling in the shared-string run, not the cell style).
    $ws.Range("E263:E263").Copy() | Out-Null
    $ws.Range("E" + $row).PasteSpecial(-4122) | Out-Null
}

Set-EventRow 264 45759 "RAVE IM REINEKE TRANCE, BOUNCE & ATZENSOUND" "Reineke Fuchs" "Köln" "https://www.instagram.com/reel/DHQftpNsoCn/?igsh=MWlqeTMxcHBxa2R6NQ=="

Set-EventRow 265 45766 "OSTER RAVE" "Artheater" "Köln" "https://www.instagram.com/s/aGlnaGxpZ2h0OjE4MDA0MzAwOTIyMTg4OTc5?story_media_id=3595802300236577132&igsh=MXg3dnIzZGE0eDkxNQ=="

Set-EventRow 266 45743 "TURBO 120 MIN RAVE" "Oma Doris" "Dortmund" "https://turbo.ticket.io/trc9thsk/"

Write-Host "Rows 264-266 populated."
